# Auto-generated Excel COM-interop script to sync the Forms export
# table "Responses" (xl/worksheets/sheet1.xml) to its refreshed state:
# corrected rater emails on several rows, a literal "utcNow()" placeholder
# written into N2, and two new response rows (9 and 10) appended, with the
# table/autofilter/used-range grown to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the "Responses" table (and its autofilter) to A1:N10 first so
# the two appended rows are already part of the table when populated. ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N10")) | Out-Null

# --- Cell values: header row is unchanged; rows 2-8 get a handful of
# corrected values; rows 9-10 are brand-new response rows. ---
# Row 1
$ws.Range("A1").Value = "ResponseId"
$ws.Range("B1").Value = "SubmittedAt"
$ws.Range("C1").Value = "RaterEmail"
$ws.Range("D1").Value = "RaterName"
$ws.Range("E1").Value = "PresenterChoice"
$ws.Range("F1").Value = "Q2_1"
$ws.Range("G1").Value = "Q2_2"
$ws.Range("H1").Value = "Q2_3"
$ws.Range("I1").Value = "Q3_1"
$ws.Range("J1").Value = "Q3_3"
$ws.Range("K1").Value = "Q3_2"
$ws.Range("L1").Value = "Q_4"
$ws.Range("M1").Value = "Comment"
$ws.Range("N1").Value = "ReceivedAtUTC"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "arda.aydin@uzh.ch"
$ws.Range("D2").Value = "arda.aydin@uzh.ch"
$ws.Range("E2").Value = "Arda`tAydin"
$ws.Range("F2").Value = " Bad 1"
$ws.Range("G2").Value = " Bad 1"
$ws.Range("H2").Value = " Bad 1"
$ws.Range("I2").Value = "Bad 1"
$ws.Range("J2").Value = "Bad 1"
$ws.Range("K2").Value = "Bad 1"
$ws.Range("L2").Value = "Bad 1"
$ws.Range("M2").Value = "testtt"
$ws.Range("N2").Value = "utcNow()"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("C3").Value = "achille.desbrieres@uzh.ch"
$ws.Range("D3").Value = "achille.desbrieres@uzh.ch"
$ws.Range("E3").Value = "Arda`tAydin"
$ws.Range("F3").Value = " Bad 1"
$ws.Range("G3").Value = " Bad 1"
$ws.Range("H3").Value = " Bad 1"
$ws.Range("I3").Value = "Bad 1"
$ws.Range("J3").Value = "Bad 1"
$ws.Range("K3").Value = "Bad 1"
$ws.Range("L3").Value = "Bad 1"
$ws.Range("M3").Value = "test2222"
$ws.Range("N3").Value = "2026-02-19T14:35:52.5219332Z"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = "achille.desbrieres@uzh.ch"
$ws.Range("D4").Value = "achille.desbrieres@uzh.ch"
$ws.Range("E4").Value = "Abirami`tSivarajah"
$ws.Range("F4").Value = "Good 4"
$ws.Range("G4").Value = "Excellent 5"
$ws.Range("H4").Value = "Fair 3"
$ws.Range("I4").Value = "Fair 3"
$ws.Range("J4").Value = "Good 4"
$ws.Range("K4").Value = "Poor 2"
$ws.Range("L4").Value = "Poor 2"
$ws.Range("M4").Value = "testcase"
$ws.Range("N4").Value = "2026-02-20T09:50:20.4657318Z"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("C5").Value = "arda.aydin@uzh.ch"
$ws.Range("D5").Value = "arda.aydin@uzh.ch"
$ws.Range("E5").Value = "Gizem`tTopsakal"
$ws.Range("F5").Value = " Bad 1"
$ws.Range("G5").Value = "Excellent 5"
$ws.Range("H5").Value = "Excellent 5"
$ws.Range("I5").Value = "Excellent 5"
$ws.Range("J5").Value = "Poor 2"
$ws.Range("K5").Value = "Bad 1"
$ws.Range("L5").Value = "Excellent 5"
$ws.Range("M5").Value = "testcase"
$ws.Range("N5").Value = "2026-02-20T09:50:50.0443991Z"

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("C6").Value = "achille.desbrieres@uzh.ch"
$ws.Range("D6").Value = "achille.desbrieres@uzh.ch"
$ws.Range("E6").Value = "Anna-Lea`tWölfle"
$ws.Range("F6").Value = "Fair 3"
$ws.Range("G6").Value = "Good 4"
$ws.Range("H6").Value = "Poor 2"
$ws.Range("I6").Value = "Poor 2"
$ws.Range("J6").Value = "Fair 3"
$ws.Range("K6").Value = "Poor 2"
$ws.Range("L6").Value = "Excellent 5"
$ws.Range("M6").Value = "testcase"
$ws.Range("N6").Value = "2026-02-20T09:51:05.1273717Z"

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("C7").Value = "achille.desbrieres@uzh.ch"
$ws.Range("D7").Value = "achille.desbrieres@uzh.ch"
$ws.Range("E7").Value = "Giacomo`tBasile"
$ws.Range("F7").Value = "Excellent 5"
$ws.Range("G7").Value = "Excellent 5"
$ws.Range("H7").Value = "Excellent 5"
$ws.Range("I7").Value = "Excellent 5"
$ws.Range("J7").Value = "Excellent 5"
$ws.Range("K7").Value = "Excellent 5"
$ws.Range("L7").Value = "Excellent 5"
$ws.Range("M7").Value = "testcase"
$ws.Range("N7").Value = "2026-02-20T09:51:20.7604177Z"

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("C8").Value = "arda.aydin@uzh.ch"
$ws.Range("D8").Value = "arda.aydin@uzh.ch"
$ws.Range("E8").Value = "Alejandro`tCaicedo Murgueitio"
$ws.Range("F8").Value = "Poor 2"
$ws.Range("G8").Value = "Fair 3"
$ws.Range("H8").Value = "Good 4"
$ws.Range("I8").Value = "Good 4"
$ws.Range("J8").Value = "Poor 2"
$ws.Range("K8").Value = "Fair 3"
$ws.Range("L8").Value = "Good 4"
$ws.Range("M8").Value = "tescase"
$ws.Range("N8").Value = "2026-02-20T09:58:26.8961547Z"

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 46075.82246527778
$ws.Range("C9").Value = "achille.desbrieres@uzh.ch"
$ws.Range("D9").Value = "achille.desbrieres@uzh.ch"
$ws.Range("E9").Value = "Steven Thomas`tUvakov"
$ws.Range("F9").Value = " Bad 1"
$ws.Range("G9").Value = " Bad 1"
$ws.Range("H9").Value = " Bad 1"
$ws.Range("I9").Value = "Bad 1"
$ws.Range("J9").Value = "Bad 1"
$ws.Range("K9").Value = "Bad 1"
$ws.Range("L9").Value = "Bad 1"
$ws.Range("M9").Value = "testttt"
$ws.Range("N9").Value = "2026-02-22T19:44:21.8423287Z"

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 46075.82340277778
$ws.Range("C10").Value = "achille.desbrieres@uzh.ch"
$ws.Range("D10").Value = "achille.desbrieres@uzh.ch"
$ws.Range("E10").Value = "Erik`tMacniel"
$ws.Range("F10").Value = "Excellent 5"
$ws.Range("G10").Value = "Excellent 5"
$ws.Range("H10").Value = "Excellent 5"
$ws.Range("I10").Value = "Excellent 5"
$ws.Range("J10").Value = "Excellent 5"
$ws.Range("K10").Value = "Excellent 5"
$ws.Range("L10").Value = "Excellent 5"
$ws.Range("M10").Value = "testestest"
$ws.Range("N10").Value = "2026-02-22T19:45:42.8842352Z"

# --- Column B (SubmittedAt) for the two new rows is a real date-time
# serial; format it with the builtin m/d/yy h:mm (numFmtId 22) style. ---
$ws.Range("B9:B10").NumberFormat = "m/d/yy h:mm"

# --- Rows whose RaterEmail/RaterName cells use the larger 12pt "Normal 2"
# font get re-themed from a theme-based black to an explicit RGB black,
# and Excel auto-fits those rows a touch taller (15.75pt) to match. ---
$boldRows = @(3,4,6,7,9,10)
$boldRange = $null
foreach ($r in $boldRows) {
    $rowRng = $ws.Range("C" + $r + ":D" + $r)
    if ($boldRange -eq $null) {
        $boldRange = $rowRng
    } else {
        $boldRange = $excel.Union($boldRange, $rowRng)
    }
    $ws.Rows.Item($r).RowHeight = 15.75
}
$boldRange.Font.Size = 12
$boldRange.Font.Color = 0

# Row 8 reverts to the plain (non-bold-font) style, and the remaining
# un-styled rows keep their default row height.
$ws.Range("C8:D8").Font.Size = 11
$ws.Range("C8:D8").Font.Color = 0

# --- Dimension / used-range grows to A1:N10 automatically once the new
# rows carry data (already satisfied above); just move the selection to
# match where the editor last left off. ---
$ws.Range("D8").Select() | Out-Null

